$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Label-column fixups -------------------------------------------------
# The commit swaps the shared strings "Inflammatory-Mac" and "MuSCs" in the
# string table. That silently relabels every cell that pointed at those two
# strings, including the "Sending cluster" values in column A for the rows
# that used to read "Inflammatory-Mac" (rows 8-13 now read "MuSCs"), and the
# "Target cluster" values in column D that explicitly flip index (rows 4/5,
# 10/11, 16/17 swap between "Inflammatory-Mac" and "MuSCs").
8..13 | ForEach-Object { $ws.Range("A$_").Value() = "MuSCs" }

$ws.Range("D4").Value() = "Inflammatory-Mac"
$ws.Range("D5").Value() = "MuSCs"
$ws.Range("D10").Value() = "Inflammatory-Mac"
$ws.Range("D11").Value() = "MuSCs"
$ws.Range("D16").Value() = "Inflammatory-Mac"
$ws.Range("D17").Value() = "MuSCs"

# --- Recomputed TPM-derived metrics (NATMI output columns E-T) -----------
# Table of (cell, new value) pulled straight from the new TPM recompute.
$edits = @(
    @("I2", 0.3637788484215491),
    @("J2", 0.3637788484215491),
    @("M2", 1.363113),
    @("N2", 4.089339),
    @("O2", 0.06061833851125786),
    @("P2", 0.06061833851125786),
    @("Q2", 0.122649727143),
    @("R2", 1.103847544287),
    @("S2", 0.02205166937685303),
    @("T2", 0.02205166937685302),
    @("I3", 0.3637788484215491),
    @("J3", 0.3637788484215491),
    @("O3", 0.1775969932713293),
    @("P3", 0.1775969932713292),
    @("S3", 0.06460602969537377),
    @("T3", 0.06460602969537375),
    @("I4", 0.3637788484215491),
    @("J4", 0.3637788484215491),
    @("M4", 3.206217333333333),
    @("N4", 9.618651999999999),
    @("O4", 0.1425821393036839),
    @("P4", 0.1425821393036839),
    @("Q4", 0.2884879544795555),
    @("R4", 2.596391590316),
    @("S4", 0.05186836644137503),
    @("T4", 0.05186836644137502),
    @("I5", 0.3637788484215491),
    @("J5", 0.3637788484215491),
    @("M5", 2.765104),
    @("N5", 8.295312000000001),
    @("O5", 0.1229656017445606),
    @("P5", 0.1229656017445605),
    @("Q5", 0.2487976060106667),
    @("R5", 2.239178454096),
    @("S5", 0.04473228499809908),
    @("T5", 0.04473228499809907),
    @("I6", 0.3637788484215491),
    @("J6", 0.3637788484215491),
    @("M6", 4.898567333333333),
    @("N6", 14.695702),
    @("O6", 0.2178418170996753),
    @("P6", 0.2178418170996753),
    @("Q6", 0.4407616586628889),
    @("R6", 3.966854927966),
    @("S6", 0.07924624536257761),
    @("T6", 0.07924624536257761),
    @("I7", 0.3637788484215491),
    @("J7", 0.3637788484215491),
    @("M7", 6.260217666666667),
    @("N7", 18.780653),
    @("O7", 0.278395110069493),
    @("P7", 0.278395110069493),
    @("Q7", 0.5632797784721111),
    @("R7", 5.069518006249),
    @("S7", 0.1012742525472706),
    @("T7", 0.1012742525472706),
    @("E8", 1),
    @("F8", 0.3333333333333333),
    @("G8", 0.039583),
    @("H8", 0.118749),
    @("I8", 0.1600336915872107),
    @("J8", 0.1600336915872107),
    @("M8", 1.363113),
    @("N8", 4.089339),
    @("O8", 0.06061833851125786),
    @("P8", 0.06061833851125786),
    @("Q8", 0.053956101879),
    @("R8", 0.4856049169109999),
    @("S8", 0.009700976489839776),
    @("T8", 0.009700976489839775),
    @("E9", 1),
    @("F9", 0.3333333333333333),
    @("G9", 0.039583),
    @("H9", 0.118749),
    @("I9", 0.1600336915872107),
    @("J9", 0.1600336915872107),
    @("O9", 0.1775969932713293),
    @("P9", 0.1775969932713292),
    @("Q9", 0.1580782597756667),
    @("R9", 1.422704337981),
    @("S9", 0.02842150244799984),
    @("T9", 0.02842150244799983),
    @("E10", 1),
    @("F10", 0.3333333333333333),
    @("G10", 0.039583),
    @("H10", 0.118749),
    @("I10", 0.1600336915872107),
    @("J10", 0.1600336915872107),
    @("M10", 3.206217333333333),
    @("N10", 9.618651999999999),
    @("O10", 0.1425821393036839),
    @("P10", 0.1425821393036839),
    @("Q10", 0.1269117007053333),
    @("R10", 1.142205306348),
    @("S10", 0.02281794610717046),
    @("T10", 0.02281794610717046),
    @("E11", 1),
    @("F11", 0.3333333333333333),
    @("G11", 0.039583),
    @("H11", 0.118749),
    @("I11", 0.1600336915872107),
    @("J11", 0.1600336915872107),
    @("M11", 2.765104),
    @("N11", 8.295312000000001),
    @("O11", 0.1229656017445606),
    @("P11", 0.1229656017445605),
    @("Q11", 0.109451111632),
    @("R11", 0.985060004688),
    @("S11", 0.01967863918542478),
    @("T11", 0.01967863918542478),
    @("E12", 1),
    @("F12", 0.3333333333333333),
    @("G12", 0.039583),
    @("H12", 0.118749),
    @("I12", 0.1600336915872107),
    @("J12", 0.1600336915872107),
    @("M12", 4.898567333333333),
    @("N12", 14.695702),
    @("O12", 0.2178418170996753),
    @("P12", 0.2178418170996753),
    @("Q12", 0.1938999907553333),
    @("R12", 1.745099916798),
    @("S12", 0.03486203017252699),
    @("T12", 0.03486203017252699),
    @("E13", 1),
    @("F13", 0.3333333333333333),
    @("G13", 0.039583),
    @("H13", 0.118749),
    @("I13", 0.1600336915872107),
    @("J13", 0.1600336915872107),
    @("M13", 6.260217666666667),
    @("N13", 18.780653),
    @("O13", 0.278395110069493),
    @("P13", 0.278395110069493),
    @("Q13", 0.2477981958996667),
    @("R13", 2.230183763097),
    @("S13", 0.04455259718424882),
    @("T13", 0.04455259718424882),
    @("E14", 2),
    @("F14", 0.6666666666666666),
    @("G14", 0.117781),
    @("H14", 0.353343),
    @("I14", 0.4761874599912402),
    @("J14", 0.4761874599912402),
    @("M14", 1.363113),
    @("N14", 4.089339),
    @("O14", 0.06061833851125786),
    @("P14", 0.06061833851125786),
    @("Q14", 0.160548812253),
    @("R14", 1.444939310277),
    @("S14", 0.02886569264456506),
    @("T14", 0.02886569264456506),
    @("E15", 2),
    @("F15", 0.6666666666666666),
    @("G15", 0.117781),
    @("H15", 0.353343),
    @("I15", 0.4761874599912402),
    @("J15", 0.4761874599912402),
    @("O15", 0.1775969932713293),
    @("P15", 0.1775969932713292),
    @("Q15", 0.4703689845296667),
    @("R15", 4.233320860767001),
    @("S15", 0.08456946112795566),
    @("T15", 0.08456946112795566),
    @("E16", 2),
    @("F16", 0.6666666666666666),
    @("G16", 0.117781),
    @("H16", 0.353343),
    @("I16", 0.4761874599912402),
    @("J16", 0.4761874599912402),
    @("M16", 3.206217333333333),
    @("N16", 9.618651999999999),
    @("O16", 0.1425821393036839),
    @("P16", 0.1425821393036839),
    @("Q16", 0.3776314837373334),
    @("R16", 3.398683353636),
    @("S16", 0.06789582675513842),
    @("T16", 0.06789582675513842),
    @("E17", 2),
    @("F17", 0.6666666666666666),
    @("G17", 0.117781),
    @("H17", 0.353343),
    @("I17", 0.4761874599912402),
    @("J17", 0.4761874599912402),
    @("M17", 2.765104),
    @("N17", 8.295312000000001),
    @("O17", 0.1229656017445606),
    @("P17", 0.1229656017445605),
    @("Q17", 0.3256767142240001),
    @("R17", 2.931090428016001),
    @("S17", 0.05855467756103671),
    @("T17", 0.0585546775610367),
    @("E18", 2),
    @("F18", 0.6666666666666666),
    @("G18", 0.117781),
    @("H18", 0.353343),
    @("I18", 0.4761874599912402),
    @("J18", 0.4761874599912402),
    @("M18", 4.898567333333333),
    @("N18", 14.695702),
    @("O18", 0.2178418170996753),
    @("P18", 0.2178418170996753),
    @("Q18", 0.5769581590873334),
    @("R18", 5.192623431786001),
    @("S18", 0.1037335415645707),
    @("T18", 0.1037335415645707),
    @("E19", 2),
    @("F19", 0.6666666666666666),
    @("G19", 0.117781),
    @("H19", 0.353343),
    @("I19", 0.4761874599912402),
    @("J19", 0.4761874599912402),
    @("M19", 6.260217666666667),
    @("N19", 18.780653),
    @("O19", 0.278395110069493),
    @("P19", 0.278395110069493),
    @("Q19", 0.7373346969976667),
    @("R19", 6.636012272979),
    @("S19", 0.1325682603379736),
    @("T19", 0.1325682603379736)
)

foreach ($edit in $edits) {
    $ws.Range($edit[0]).Value() = $edit[1]
}
